$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.123.82'
$ws.Range('E2').Value = '  +0.35%  '
$ws.Range('D3').Value = '2.685.65'
$ws.Range('E3').Value = '  +4.89%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '''518.63'
$ws.Range('E5').Value = '  +2.52%  '
$ws.Range('D6').Value = '''146.12'
$ws.Range('E6').Value = '  +2.84%  '
$ws.Range('D7').Value = '''0.995'
$ws.Range('E7').Value = '  -0.39%  '
$ws.Range('D8').Value = '''0.569'
$ws.Range('E8').Value = '  +2.90%  '
$ws.Range('D9').Value = '2.719.60'
$ws.Range('E9').Value = '  +6.16%  '
$ws.Range('D10').Value = '''6.27'
$ws.Range('E10').Value = '  +1.06%  '
$ws.Range('E11').Value = '  +6.87%  '
$ws.Range('E12').Value = '  +2.84%  '
$ws.Range('E13').Value = '  -0.94%  '
$ws.Range('D14').Value = '3.156.64'
$ws.Range('E14').Value = '  +5.00%  '
$ws.Range('D15').Value = '59.110.41'
$ws.Range('E15').Value = '  +0.37%  '
$ws.Range('D16').Value = '''21.20'
$ws.Range('E16').Value = '  +3.42%  '
$ws.Range('E17').Value = '  +3.30%  '
$ws.Range('D18').Value = '2.720.99'
$ws.Range('E18').Value = '  +6.18%  '
$ws.Range('D19').Value = '''4.57'
$ws.Range('E19').Value = '  +1.70%  '
$ws.Range('E20').Value = '  +4.80%  '
$ws.Range('D21').Value = '''10.52'
$ws.Range('E21').Value = '  +4.92%  '
$ws.Range('D22').Value = '''6.23'
$ws.Range('E22').Value = '  +5.24%  '
$ws.Range('E23').Value = '  -0.35%  '
$ws.Range('D24').Value = '''61.15'
$ws.Range('E24').Value = '  +2.84%  '
$ws.Range('D25').Value = '''0.424'
$ws.Range('E25').Value = '  +5.01%  '
$ws.Range('D26').Value = '2.798.15'
$ws.Range('E26').Value = '  +4.86%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').Value = '''0.991'
$ws.Range('E27').Value = '  -0.75%  '
$ws.Range('B28').Value = 'Kaspa'
$ws.Range('C28').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D28').Value = '''0.161'
$ws.Range('E28').Value = '  +2.47%  '
$ws.Range('D29').Value = '0.0₃0826'
$ws.Range('E29').Value = '  +6.75%  '
$ws.Range('E30').Value = '  +6.30%  '
$ws.Range('D31').Value = '''0.997'
$ws.Range('E31').Value = '  -0.32%  '
$ws.Range('D32').Value = '''6.45'
$ws.Range('E32').Value = '  +11.60%  '
$ws.Range('D33').Value = '''19.16'
$ws.Range('E33').Value = '  +3.65%  '
$ws.Range('E34').Value = '  +3.44%  '
$ws.Range('D35').Value = '''150.20'
$ws.Range('E35').Value = '  +0.69%  '
$ws.Range('E36').Value = '  +17.58%  '
$ws.Range('D37').Value = '''4.08'
$ws.Range('E37').Value = '  +5.47%  '
$ws.Range('E38').Value = '  +5.31%  '
$ws.Range('E39').Value = '  +5.06%  '
$ws.Range('D40').Value = '''36.90'
$ws.Range('E40').Value = '  +3.02%  '
$ws.Range('E41').Value = '  +6.97%  '
$ws.Range('E42').Value = '  +3.04%  '
$ws.Range('D43').Value = '''0.627'
$ws.Range('E43').Value = '  +3.67%  '
$ws.Range('D44').Value = '''282.70'
$ws.Range('E44').Value = '  -1.13%  '
$ws.Range('D45').Value = '''20.32'
$ws.Range('E45').Value = '  +9.50%  '
$ws.Range('E46').Value = '  +0.63%  '
$ws.Range('D47').Value = '''0.992'
$ws.Range('E47').Value = '  -0.58%  '
$ws.Range('D48').Value = '''0.0537'
$ws.Range('E48').Value = '  +1.68%  '
$ws.Range('D49').Value = '''4.79'
$ws.Range('E49').Value = '  +6.24%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').Value = '''0.0232'
$ws.Range('E50').Value = '  +3.20%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '2.021.34'
$ws.Range('E51').Value = '  +7.36%  '
